$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.771.40"
$ws.Range("E2").Value = "  +4.73%  "

$ws.Range("D3").Value = "2.280.14"
$ws.Range("E3").Value = "  +2.25%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "231.24"
$ws.Range("E5").Value = "  -0.36%  "

$ws.Range("E6").Value = "  +0.74%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "64.78"
$ws.Range("E7").Value = "  +7.12%  "

$ws.Range("E8").Value = "  -0.02%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.426"
$ws.Range("E9").Value = "  +5.19%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0969"
$ws.Range("E10").Value = "  +7.68%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "57.88"
$ws.Range("E11").Value = "  -0.44%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "26.41"
$ws.Range("E12").Value = "  +16.48%  "

$ws.Range("E13").Value = "  +0.43%  "

$ws.Range("D14").Value = "2.618.21"
$ws.Range("E14").Value = "  +2.26%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.79"
$ws.Range("E15").Value = "  +1.29%  "

$ws.Range("E16").Value = "  +5.33%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.820"
$ws.Range("E17").Value = "  +2.32%  "

$ws.Range("D18").Value = "2.304.28"
$ws.Range("E18").Value = "  +3.36%  "

$ws.Range("D19").Value = "43.596.48"
$ws.Range("E19").Value = "  +4.50%  "

$ws.Range("E20").Value = "  +4.58%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "73.38"
$ws.Range("E21").Value = "  +1.18%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.17"
$ws.Range("E22").Value = "  +0.03%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "250.42"
$ws.Range("E23").Value = "  +0.95%  "

$ws.Range("E24").Value = "  +0.11%  "

$ws.Range("E25").Value = "  +5.80%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.46"
$ws.Range("E26").Value = "  +3.21%  "

$ws.Range("E27").Value = "  +3.53%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "172.12"
$ws.Range("E28").Value = "  +1.50%  "

$ws.Range("E29").Value = "  -2.68%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "20.58"
$ws.Range("E30").Value = "  +3.25%  "

$ws.Range("E31").Value = "  +3.03%  "

$ws.Range("E33").Value = "  +0.27%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.24"
$ws.Range("E34").Value = "  +4.48%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0698"
$ws.Range("E35").Value = "  +6.66%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.73"
$ws.Range("E36").Value = "  +0.65%  "

$ws.Range("E37").Value = "  +3.60%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.79"
$ws.Range("E38").Value = "  +4.64%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.36"
$ws.Range("E39").Value = "  -1.12%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0247"
$ws.Range("E40").Value = "  +3.07%  "

$ws.Range("E41").Value = "  -0.08%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "11.03"
$ws.Range("E42").Value = "  +28.04%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.000228"
$ws.Range("E43").Value = "  -3.37%  "

$ws.Range("E44").Value = "  +4.11%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.47"
$ws.Range("E45").Value = "  -1.50%  "

$ws.Range("E46").Value = "  -0.10%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0966"
$ws.Range("E47").Value = "  +0.33%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "97.98"
$ws.Range("E48").Value = "  -0.80%  "

$ws.Range("D49").Value = "1.487.77"
$ws.Range("E49").Value = "  +1.15%  "

$ws.Range("E50").Value = "  +1.69%  "

$ws.Range("E51").Value = "  +0.87%  "
